$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.707.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.055.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '513.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.21%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.33'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.106'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.375'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.578.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.71%  '
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.730.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.060.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("E18").Value = '  +4.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.497'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  +3.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0894'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.79'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '155.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0673'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.095.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.651'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.269.92'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0253'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.16%  '
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.929'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.723'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '254.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.62%  '
